$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "active" column header in K7, matching the style used by the
# other header cells in row 7 (copy format from J7, then set the value).
$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("K7").Value = "active"

# Update the selection to match the new target range for the report body.
[void]$ws.Range("A8:A41").Select()
